$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 45 data (set title/author first so shared-string order matches)
$ws.Range("A45").Value = "Scene Memory Transformer for Embodied Agents in Long-Horizon Tasks"
$ws.Range("B45").Value = 2019
$ws.Range("C45").Value = "Kuan Fang, Alexander Toshev, Li Fei-Fei, Silvio Savarese"
$ws.Range("D45").Value = "seqence to seqence"
$ws.Range("E45").Value = "attention model"

# Update G37: "powered by Google" -> "Transformer powered by Google"
# (new shared string gets created here, then reused by G45)
$ws.Range("G37").Value = "Transformer powered by Google"
$ws.Range("G45").Value = "Transformer powered by Google"

# Add F37 date (2020-05-19, serial 43970) with date formatting like F41/F4
$ws.Range("F37").Value = 43970
$ws.Range("F41").Copy()
$ws.Range("F37").PasteSpecial(-4122)

# Resize the table (ListObject) to include the new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G45"))

# Update sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$ws.Range("F47").Select()
